# Add nowcasts for 2025Q4
#
# The sheet holds a 7x11 "nowcast" table: a header row (Row, Prognose,
# surveys, production, orders, turnover, financial, labor market, prices,
# national accounts, Revision) followed by 6 data rows keyed by a
# "vintage" date in column A. This update rolls the 6 vintage dates
# forward by one quarter (2025-06-30..2025-09-15 -> 2025-09-30..2025-12-15)
# and refreshes the corresponding nowcast/revision figures in B:K.
# Rows 8-11 (the next block of vintages) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Header row (unchanged text, rewritten for completeness) --
$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "Prognose"
$ws.Range("C1").Value = "surveys"
$ws.Range("D1").Value = "production"
$ws.Range("E1").Value = "orders"
$ws.Range("F1").Value = "turnover"
$ws.Range("G1").Value = "financial"
$ws.Range("H1").Value = "labor market"
$ws.Range("I1").Value = "prices"
$ws.Range("J1").Value = "national accounts"
$ws.Range("K1").Value = "Revision"

# -- Column A holds plain text vintage dates (not real Excel dates), so
#    force a text format before writing them; this keeps the cells typed
#    as strings instead of being auto-converted to date serials.
$ws.Range("A2:A7").NumberFormat = "@"

# -- Row 2 : vintage 2025-09-30 --
$ws.Range("A2").Value = "2025-09-30"
$ws.Range("B2").Value = 0.26892855483751477
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

# -- Row 3 : vintage 2025-10-15 --
$ws.Range("A3").Value = "2025-10-15"
$ws.Range("B3").Value = 0.11295526046108378
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -0.19715656854375554
$ws.Range("E3").Value = -0.0029555110305725316
$ws.Range("F3").Value = 0.016831949170904529
$ws.Range("G3").Value = 0.0068804160627636566
$ws.Range("H3").Value = 0.002238106331374
$ws.Range("I3").Value = 0.018906192219394657
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = -0.00071787858653971082

# -- Row 4 : vintage 2025-10-30 --
$ws.Range("A4").Value = "2025-10-30"
$ws.Range("B4").Value = 0.5114491829359662
$ws.Range("C4").Value = 0.16394768541553725
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.00061000002029408431
$ws.Range("F4").Value = 0.0002011988465820333
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = -0.01048990809688379
$ws.Range("I4").Value = 0.22111278609876595
$ws.Range("J4").Value = 0.020945074396409562
$ws.Range("K4").Value = 0.002167085794177348

# -- Row 5 : vintage 2025-11-15 --
$ws.Range("A5").Value = "2025-11-15"
$ws.Range("B5").Value = 0.33147271249921606
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = -0.073833175829404726
$ws.Range("E5").Value = 0.012177545168242718
$ws.Range("F5").Value = -0.15656258394572245
$ws.Range("G5").Value = -0.0061390117967748629
$ws.Range("H5").Value = -0.0091858915794224109
$ws.Range("I5").Value = 0.021555375296689251
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0.032011272249642397

# -- Row 6 : vintage 2025-11-30 --
$ws.Range("A6").Value = "2025-11-30"
$ws.Range("B6").Value = 0.028256432548483923
$ws.Range("C6").Value = -0.20797927517932896
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.0016282621350707264
$ws.Range("F6").Value = -0.0049240913733429965
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = -0.0086845894728944922
$ws.Range("I6").Value = -0.084736416405899406
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0.0014798303456630135

# -- Row 7 : vintage 2025-12-15 --
$ws.Range("A7").Value = "2025-12-15"
$ws.Range("B7").Value = 0.10583447526029355
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.22631323161650063
$ws.Range("E7").Value = 0.00073343107051046974
$ws.Range("F7").Value = -0.12360441206627915
$ws.Range("G7").Value = 0.004130969187833301
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = -0.029995177096755643

# Drop back to the default (general) style now that the text is in place,
# so the cells don't carry a stray text-format style index.
$ws.Range("A2:A7").Style = "Normal"

# -- Column widths: nudge to the widths that accompanied the new content --
$ws.Columns.Item(2).ColumnWidth = 12.833333333333332
$ws.Columns.Item(3).ColumnWidth = 12.333333333333332
$ws.Columns.Item(4).ColumnWidth = 13.333333333333332
$ws.Columns.Item(5).ColumnWidth = 14.833333333333332
$ws.Columns.Item(6).ColumnWidth = 14.833333333333332
$ws.Columns.Item(7).ColumnWidth = 14.333333333333332
$ws.Columns.Item(8).ColumnWidth = 14.333333333333332
$ws.Columns.Item(9).ColumnWidth = 13.333333333333332
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
